$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update sheet1 ("o_10") ---

# Add new header column E ("evaluator_partial_correctness"), copying the
# formatting from the existing D1 header cell so the style matches.
$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "evaluator_partial_correctness"

# New data for row 2.
$prompt = " Given is the adjacency matrix for a unweighted undirected graph containing 14 nodes labelled A to N. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   `n`nConsider some examples`n`nExample 1: what is the shortest path from node A to node G?`n   A B C D E F G H I`n A 0 1 0 0 0 0 0 1 0`n B 1 0 1 0 0 1 0 1 1`n C 0 1 0 1 0 0 0 0 0`n D 0 0 1 0 1 0 0 0 0`n E 0 0 0 1 0 0 0 0 0`n F 0 1 0 0 0 0 1 0 1`n G 0 0 0 0 0 1 0 0 0`n H 1 1 0 0 0 0 0 0 0`n I 0 1 0 0 0 1 0 0 0`n`nSolution: A -> B -> F -> G`n        `n`nExample 2: what is the shortest path from node A to node G?`n   A B C D E F G H I J K`n A 0 1 1 0 0 0 0 1 0 0 0`n B 1 0 0 0 0 0 0 0 0 0 0`n C 1 0 0 1 1 0 0 0 0 0 0`n D 0 0 1 0 0 0 0 0 0 0 0`n E 0 0 1 0 0 1 1 0 0 0 1`n F 0 0 0 0 1 0 0 0 0 0 0`n G 0 0 0 0 1 0 0 0 0 1 1`n H 1 0 0 0 0 0 0 0 1 0 0`n I 0 0 0 0 0 0 0 1 0 1 0`n J 0 0 0 0 0 0 1 0 1 0 0`n K 0 0 0 0 1 0 1 0 0 0 0`n`nSolution: A -> C -> E -> G`n        `n`nExample 3: what is the shortest path from node A to node O?`n   A B C D E F G H I J K L M N O P Q R S T U`n A 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0`n B 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0`n C 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0`n D 0 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0`n E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0`n F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0`n G 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0`n H 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0`n I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 1 0 0`n J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0`n K 0 0 0 0 0 0 0 0 0 1 0 1 1 1 0 0 0 0 0 0 1`n L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0`n M 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0`n N 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0`n O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0`n P 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0`n Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0`n R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0`n S 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 0`n T 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1`n U 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0`n`nSolution: A -> T -> U -> K -> N -> O`n        `n Given these examples, answer the following quesiton.`n`nwhat is the shortest path from node A to node J?`n`n   A B C D E F G H I J K L M N`n A 0 1 0 0 0 0 1 0 0 0 0 0 0 0`n B 1 0 1 0 1 0 0 0 0 0 0 0 0 0`n C 0 1 0 1 0 0 0 0 0 0 0 0 0 0`n D 0 0 1 0 0 0 0 0 0 0 0 0 0 0`n E 0 1 0 0 0 1 0 0 0 0 0 0 0 0`n F 0 0 0 0 1 0 0 0 0 0 0 0 0 0`n G 1 0 0 0 0 0 0 1 0 0 0 0 0 0`n H 0 0 0 0 0 0 1 0 1 1 0 0 0 1`n I 0 0 0 0 0 0 0 1 0 0 0 0 0 0`n J 0 0 0 0 0 0 0 1 0 0 1 0 0 0`n K 0 0 0 0 0 0 0 0 0 1 0 1 0 0`n L 0 0 0 0 0 0 0 0 0 0 1 0 1 0`n M 0 0 0 0 0 0 0 0 0 0 0 1 0 1`n N 0 0 0 0 0 0 0 1 0 0 0 0 1 0`n    "
$solution = "A -> G -> H -> J"
$llm_response = "The shortest path from node A to node J in the given graph can be found by using a search algorithm such as BFS (breadth-first search) or Dijkstra's algorithm.`nUsing BFS, we can start by visiting node A and exploring its adjacent nodes. We'll keep track of the visited nodes and their distance from the starting node.`nStarting with node A:`n- Visit node B and mark it as visited with distance 1.`n- Visit node G and mark it as visited with distance 1.`n- Visit node H and mark it as visited with distance 1.`n  - Visit node I and mark it as visited with distance 2.`n    - Visit node J and mark it as visited with distance 3.`nSo, the shortest path from node A to node J is A -> H -> I -> J. The distance of this path is 3."
$evaluator_response = "Wrong"
$evaluator_partial_correctness = "Output: 2/4"

$ws1.Range("A2").Value = $prompt
$ws1.Range("B2").Value = $solution
$ws1.Range("C2").Value = $llm_response
$ws1.Range("D2").Value = $evaluator_response
$ws1.Range("E2").Value = $evaluator_partial_correctness

# Undo the automatic row-height expansion triggered by the long, multi-line
# prompt/response text so row 2 keeps using the sheet's default height.
$ws1.Rows.Item(2).AutoFit()

# --- Add the two new worksheets, copying the header row format/values ---

$newSheet1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet1.Name = "o_20"
$ws1.Range("A1:E1").Copy($newSheet1.Range("A1"))

$newSheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $newSheet1)
$newSheet2.Name = "o_20_jumbled"
$ws1.Range("A1:E1").Copy($newSheet2.Range("A1"))

# Restore sheet1 as the active/selected sheet.
$ws1.Activate()
$ws1.Range("A1").Select()
